$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "worst-fit-algorithm"
$ws.Range("C2").Value = 317
$ws.Range("D2").Value = 664
$ws.Range("E2").Value = 0.4774096385542169
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = 0.3
$ws.Range("I2").Value = 4250
$ws.Range("J2").Value = 3586
$ws.Range("K2").Value = 664
$ws.Range("L2").Value = 20
$ws.Range("M2").Value = 9
